$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 390.31033
$ws.Cells.Item(33, 9).Value = 347.21738
$ws.Cells.Item(33, 10).Value = 555.5
$ws.Cells.Item(33, 11).Value = 347.21738
$ws.Cells.Item(33, 12).Value = 555.5
$ws.Cells.Item(33, 13).Value = -118.21738
$ws.Cells.Item(33, 14).Value = -1013.5

# ALC row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).ClearContents()
$ws.Cells.Item(68, 14).Value = 0

# ALC row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).ClearContents()
$ws.Cells.Item(71, 14).Value = 0

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 488.95
$ws.Cells.Item(80, 9).Value = 354
$ws.Cells.Item(80, 10).Value = 533.93335
$ws.Cells.Item(80, 11).Value = 1062
$ws.Cells.Item(80, 12).Value = 1601.80005
$ws.Cells.Item(80, 13).Value = -64
$ws.Cells.Item(80, 14).Value = -3597.80005

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 488.95
$ws.Cells.Item(83, 9).Value = 354
$ws.Cells.Item(83, 10).Value = 533.93335
$ws.Cells.Item(83, 11).Value = 3186
$ws.Cells.Item(83, 12).Value = 4805.40015
$ws.Cells.Item(83, 13).Value = 1806
$ws.Cells.Item(83, 14).Value = -14789.40015

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 2061612.4
$ws.Cells.Item(88, 9).Value = 1995
$ws.Cells.Item(88, 10).Value = 2473535.8
$ws.Cells.Item(88, 11).Value = 1995
$ws.Cells.Item(88, 12).Value = 2473535.8
$ws.Cells.Item(88, 13).Value = -1589
$ws.Cells.Item(88, 14).Value = -2474347.8

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 2061612.4
$ws.Cells.Item(91, 9).Value = 1995
$ws.Cells.Item(91, 10).Value = 2473535.8
$ws.Cells.Item(91, 11).Value = 1995
$ws.Cells.Item(91, 12).Value = 2473535.8
$ws.Cells.Item(91, 13).Value = -591
$ws.Cells.Item(91, 14).Value = -2476343.8

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3530.9187
$ws.Cells.Item(32, 9).Value = 3393.8555
$ws.Cells.Item(32, 10).Value = 7323
$ws.Cells.Item(32, 11).Value = 3393.8555
$ws.Cells.Item(32, 12).Value = 7323
$ws.Cells.Item(32, 13).Value = -3106.8555
$ws.Cells.Item(32, 14).Value = -7897

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1289.2051
$ws.Cells.Item(122, 9).Value = 1036.1613
$ws.Cells.Item(122, 11).Value = 3108.4839
$ws.Cells.Item(122, 13).Value = -658.4839000000002

# BSM row 69
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(69, 8).Value = 16500
$ws.Cells.Item(69, 10).Value = 16500
$ws.Cells.Item(69, 12).Value = 16500
$ws.Cells.Item(69, 14).Value = -18122

# BSM row 72
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(72, 8).Value = 16500
$ws.Cells.Item(72, 10).Value = 16500
$ws.Cells.Item(72, 12).Value = 49500
$ws.Cells.Item(72, 14).Value = -57612

# CRP row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).ClearContents()
$ws.Cells.Item(64, 14).Value = 0

# CRP row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).ClearContents()
$ws.Cells.Item(67, 14).Value = 0

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 4478819
$ws.Cells.Item(86, 9).Value = 13369533
$ws.Cells.Item(86, 10).Value = 33462.4
$ws.Cells.Item(86, 11).Value = 13369533
$ws.Cells.Item(86, 12).Value = 33462.4
$ws.Cells.Item(86, 13).Value = -13368410
$ws.Cells.Item(86, 14).Value = -35708.4

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 4478819
$ws.Cells.Item(89, 9).Value = 13369533
$ws.Cells.Item(89, 10).Value = 33462.4
$ws.Cells.Item(89, 11).Value = 66847665
$ws.Cells.Item(89, 12).Value = 167312
$ws.Cells.Item(89, 13).Value = -66842049
$ws.Cells.Item(89, 14).Value = -178544

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4389.5
$ws.Cells.Item(122, 9).Value = 4782.72
$ws.Cells.Item(122, 11).Value = 14348.16
$ws.Cells.Item(122, 13).Value = -11898.16

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 17243336
$ws.Cells.Item(134, 9).Value = 2018.5416
$ws.Cells.Item(134, 11).Value = 6055.6248
$ws.Cells.Item(134, 13).Value = -3520.6248

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 6689
$ws.Cells.Item(56, 9).Value = 6689
$ws.Cells.Item(56, 11).Value = 6689
$ws.Cells.Item(56, 13).Value = -6159

# CUL row 96
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(96, 8).Value = 8488.888999999999
$ws.Cells.Item(96, 10).Value = 8488.888999999999
$ws.Cells.Item(96, 12).Value = 25466.667
$ws.Cells.Item(96, 14).Value = -29584.667

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 11905875
$ws.Cells.Item(131, 9).Value = 250000300
$ws.Cells.Item(131, 10).Value = 1153.2125
$ws.Cells.Item(131, 11).Value = 750000900
$ws.Cells.Item(131, 12).Value = 3459.6375
$ws.Cells.Item(131, 13).Value = -749995860
$ws.Cells.Item(131, 14).Value = -13539.6375

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 1760.3513
$ws.Cells.Item(139, 9).Value = 1797.7391
$ws.Cells.Item(139, 10).Value = 1698.9286
$ws.Cells.Item(139, 11).Value = 5393.2173
$ws.Cells.Item(139, 12).Value = 5096.7858
$ws.Cells.Item(139, 13).Value = -253.2173000000003
$ws.Cells.Item(139, 14).Value = -15376.7858

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7257.143
$ws.Cells.Item(80, 9).Value = 7866.6665
$ws.Cells.Item(80, 11).Value = 7866.6665
$ws.Cells.Item(80, 13).Value = -6868.6665

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 7257.143
$ws.Cells.Item(83, 9).Value = 7866.6665
$ws.Cells.Item(83, 11).Value = 39333.3325
$ws.Cells.Item(83, 13).Value = -34341.3325

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(134, 8).Value = 24031.285
$ws.Cells.Item(134, 10).Value = 24031.285
$ws.Cells.Item(134, 12).Value = 72093.855
$ws.Cells.Item(134, 14).Value = -77163.855

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1901.8
$ws.Cells.Item(7, 9).Value = 1876
$ws.Cells.Item(7, 10).Value = 2005
$ws.Cells.Item(7, 11).Value = 1876
$ws.Cells.Item(7, 12).Value = 2005
$ws.Cells.Item(7, 13).Value = -1764
$ws.Cells.Item(7, 14).Value = -2229

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7070
$ws.Cells.Item(40, 9).Value = 3076.25
$ws.Cells.Item(40, 11).Value = 3076.25
$ws.Cells.Item(40, 13).Value = -2940.25

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2132.2222
$ws.Cells.Item(68, 9).Value = 2098.75
$ws.Cells.Item(68, 11).Value = 2098.75
$ws.Cells.Item(68, 13).Value = -1349.75

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 2132.2222
$ws.Cells.Item(71, 9).Value = 2098.75
$ws.Cells.Item(71, 11).Value = 10493.75
$ws.Cells.Item(71, 13).Value = -6749.75

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 1901.8
$ws.Cells.Item(126, 9).Value = 1876
$ws.Cells.Item(126, 10).Value = 2005
$ws.Cells.Item(126, 11).Value = 5628
$ws.Cells.Item(126, 12).Value = 6015
$ws.Cells.Item(126, 13).Value = -3158
$ws.Cells.Item(126, 14).Value = -10955

# WVR row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 13).ClearContents()

# WVR row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).ClearContents()
$ws.Cells.Item(119, 14).Value = 0
